$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "phone" column of aadhaar-linked numbers, loaded for the DB import
$ws.Range("B2").Value = 4154654879
$ws.Range("B3").Value = 4167329486
$ws.Range("B4").Value = 3068815942
$ws.Range("B5").Value = 4025176530

# Best-fit the new column to its content (mirrors the double-click-border autofit)
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).ColumnWidth = 10.1666666666667

# Leave the cursor where data entry finished
$ws.Range("B5").Select() | Out-Null
